$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ConfigDict")

$ws.Range("B3").Value = "El cliente ha estado 5 segundos en Hold, por favor retoma la llamada cuanto antes"
$ws.Range("B5").Value = "El cliente ha estado 10 segundos en Hold. Cuando aceptes este mensaje la llamada se retomara."

$ws.Range("B5").Select()
